# Updated cryptos list on Wed Jan 31 09:35:34 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value would otherwise be mis-parsed as a
# number (single decimal point) by plain .Value assignment. Force text
# storage so they match the original inline-string/shared-string type,
# then restore the default "Normal" style so no stray number-format
# style sticks around on the cell (matches the unstyled source cells).
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "42.799.50"
$ws.Range("E2").Value = "  -1.36%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.323.34"
$ws.Range("E3").Value = "  +0.74%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 (BNB)
Set-TextValue "D5" "304.98"
$ws.Range("E5").Value = "  -1.87%  "

# Row 6 (Solana)
Set-TextValue "D6" "100.27"
$ws.Range("E6").Value = "  -2.92%  "

# Row 7 (XRP)
$ws.Range("E7").Value = "  -4.84%  "

# Row 8 (USDC)
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 (Cardano)
Set-TextValue "D9" "0.503"
$ws.Range("E9").Value = "  -4.67%  "

# Row 10 (Avalanche)
Set-TextValue "D10" "34.36"
$ws.Range("E10").Value = "  -6.04%  "

# Row 11 (OKB)
Set-TextValue "D11" "52.07"
$ws.Range("E11").Value = "  -0.17%  "

# Row 12 (Dogecoin)
Set-TextValue "D12" "0.0792"
$ws.Range("E12").Value = "  -2.39%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  +0.66%  "

# Row 14 (Polkadot)
Set-TextValue "D14" "6.73"
$ws.Range("E14").Value = "  -4.26%  "

# Row 15 (Chainlink)
Set-TextValue "D15" "15.66"
$ws.Range("E15").Value = "  +3.87%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "2.330.17"
$ws.Range("E16").Value = "  +1.27%  "

# Row 17 (Polygon)
Set-TextValue "D17" "0.819"
$ws.Range("E17").Value = "  +1.15%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "42.722.35"
$ws.Range("E18").Value = "  -1.29%  "

# Row 19: swap - was ShibaInu, now Uniswap
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D19" "6.13"
$ws.Range("E19").Value = "  -0.78%  "

# Row 20: swap - was Uniswap, now ShibaInu
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  -2.66%  "

# Row 21 (InternetComputer(DFINITY))
Set-TextValue "D21" "11.54"
$ws.Range("E21").Value = "  -5.51%  "

# Row 22 (Litecoin)
Set-TextValue "D22" "69.02"
$ws.Range("E22").Value = "  +1.28%  "

# Row 23 (BitcoinCash)
Set-TextValue "D23" "235.19"
$ws.Range("E23").Value = "  -3.33%  "

# Row 24 (ImmutableX)
Set-TextValue "D24" "1.97"
$ws.Range("E24").Value = "  -2.42%  "

# Row 25 (PancakeSwap)
$ws.Range("E25").Value = "  -3.11%  "

# Row 26 (Dai)
$ws.Range("E26").Value = "  +0.24%  "

# Row 27: swap - was EthereumClassic, now LEO
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D27" "3.96"
$ws.Range("E27").Value = "  -0.56%  "

# Row 28: swap - was LEO, now EthereumClassic
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "25.25"
$ws.Range("E28").Value = "  +1.77%  "

# Row 29 (Toncoin)
Set-TextValue "D29" "2.18"
$ws.Range("E29").Value = "  -4.99%  "

# Row 30 (InjectiveProtocol)
Set-TextValue "D30" "34.68"
$ws.Range("E30").Value = "  -6.26%  "

# Row 31 (Cosmos)
Set-TextValue "D31" "9.19"
$ws.Range("E31").Value = "  -4.79%  "

# Row 32 (Monero)
Set-TextValue "D32" "160.57"
$ws.Range("E32").Value = "  -4.31%  "

# Row 33 (FirstDigitalUSD)
$ws.Range("E33").Value = "  +0.06%  "

# Row 34 (Filecoin)
Set-TextValue "D34" "5.04"
$ws.Range("E34").Value = "  -4.59%  "

# Row 35: swap - was RenderToken, now WEMIXToken
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D35" "2.45"
$ws.Range("E35").Value = "  -3.25%  "

# Row 36: swap - was WEMIXToken, now RenderToken
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D36" "4.56"
$ws.Range("E36").Value = "  +2.29%  "

# Row 37 (Hedera)
Set-TextValue "D37" "0.0716"
$ws.Range("E37").Value = "  -3.81%  "

# Row 38 (Celestia)
Set-TextValue "D38" "16.97"
$ws.Range("E38").Value = "  -6.66%  "

# Row 39 (LidoDAOToken)
Set-TextValue "D39" "2.87"
$ws.Range("E39").Value = "  -5.85%  "

# Row 40 (ARBITRUM)
Set-TextValue "D40" "1.82"
$ws.Range("E40").Value = "  -3.02%  "

# Row 41 (Kaspa)
$ws.Range("E41").Value = "  -4.88%  "

# Row 42 (Stellar)
$ws.Range("E42").Value = "  -3.48%  "

# Row 43 (ApeXProtocol)
$ws.Range("E43").Value = "  -6.82%  "

# Row 44 (Maker)
$ws.Range("D44").Value = "2.004.55"
$ws.Range("E44").Value = "  +1.25%  "

# Row 45 (VeChain)
Set-TextValue "D45" "0.0281"
$ws.Range("E45").Value = "  -4.49%  "

# Row 46 (EnergySwap)
Set-TextValue "D46" "18.54"
$ws.Range("E46").Value = "  -2.33%  "

# Row 47 (FraxShare)
Set-TextValue "D47" "10.14"

# Row 48 (NEARProtocol)
Set-TextValue "D48" "2.87"
$ws.Range("E48").Value = "  -4.50%  "

# Row 49 (MultiversX)
Set-TextValue "D49" "55.20"
$ws.Range("E49").Value = "  -1.43%  "

# Row 50 (HuobiToken)
Set-TextValue "D50" "2.88"
$ws.Range("E50").Value = "  -2.10%  "

# Row 51 (RocketPoolETH)
$ws.Range("D51").Value = "2.552.49"
$ws.Range("E51").Value = "  +0.81%  "
